# Applies updated profit-tracking figures (H/I/J/K/L/M/N columns) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 1043
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# row 31
$ws.Range("H31").Value = 56.8
$ws.Range("I31").Value = 56.8
$ws.Range("K31").Value = 170.4
$ws.Range("M31").Value = 59.60000000000002
# row 43
$ws.Range("H43").Value = 6857.5713
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4931
# row 132
$ws.Range("H132").Value = 2319
$ws.Range("I132").Value = 1366.1333
$ws.Range("K132").Value = 4098.3999
$ws.Range("M132").Value = -1568.3999
# row 136
$ws.Range("H136").Value = 49999
$ws.Range("J136").Value = 49999
$ws.Range("L136").Value = 49999
$ws.Range("N136").Value = -60199
# row 138
$ws.Range("H138").Value = 5266.684
$ws.Range("J138").Value = 5537.409
$ws.Range("L138").Value = 16612.227
$ws.Range("N138").Value = -26892.227
# row 140
$ws.Range("H140").Value = 89999.8
$ws.Range("J140").Value = 89999.8
$ws.Range("L140").Value = 89999.8
$ws.Range("N140").Value = -100359.8

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3645.7646
$ws.Range("I32").Value = 2207.125
$ws.Range("K32").Value = 2207.125
$ws.Range("M32").Value = -1920.125
# row 61
$ws.Range("H61").Value = 3266.6667
$ws.Range("I61").Value = 3150
$ws.Range("K61").Value = 3150
$ws.Range("M61").Value = -2938
# row 136
$ws.Range("H136").Value = 3266.6667
$ws.Range("I136").Value = 3150
$ws.Range("K136").Value = 9450
$ws.Range("M136").Value = -6900
# row 139
$ws.Range("H139").Value = 77779.5
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 2389.0952
$ws.Range("I86").Value = 2473.111
$ws.Range("K86").Value = 2473.111
$ws.Range("M86").Value = -1350.111
# row 89
$ws.Range("H89").Value = 2389.0952
$ws.Range("I89").Value = 2473.111
$ws.Range("K89").Value = 12365.555
$ws.Range("M89").Value = -6749.555
# row 134
$ws.Range("H134").Value = 1672.9445
$ws.Range("I134").Value = 1477.2354
$ws.Range("K134").Value = 4431.706200000001
$ws.Range("M134").Value = -1896.706200000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 4
$ws.Range("H4").Value = 9477.888999999999
$ws.Range("I4").Value = 2300.5
$ws.Range("J4").Value = 11528.571
$ws.Range("K4").Value = 2300.5
$ws.Range("L4").Value = 11528.571
$ws.Range("M4").Value = -2188.5
$ws.Range("N4").Value = -11752.571

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
# row 68
$ws.Range("H68").Value = 998.5
$ws.Range("I68").Value = 998
$ws.Range("K68").Value = 2994
$ws.Range("M68").Value = -2183
# row 71
$ws.Range("H71").Value = 998.5
$ws.Range("I71").Value = 998
$ws.Range("K71").Value = 8982
$ws.Range("M71").Value = -4926
# row 114
$ws.Range("H114").Value = 2477.5
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
# row 129
$ws.Range("H129").Value = 4248.875
$ws.Range("J129").Value = 4498.5
$ws.Range("L129").Value = 13495.5
$ws.Range("N129").Value = -23495.5
# row 137
$ws.Range("H137").Value = 2836.7144
$ws.Range("J137").Value = 4206.75
$ws.Range("L137").Value = 12620.25
$ws.Range("N137").Value = -22820.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 306.05884
$ws.Range("I2").Value = 10.9
$ws.Range("K2").Value = 10.9
$ws.Range("M2").Value = 102.1
# row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# row 141
$ws.Range("H141").Value = 48000
$ws.Range("J141").Value = 48000
$ws.Range("L141").Value = 48000
$ws.Range("N141").Value = -58360

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 43
$ws.Range("H43").Value = 12500
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12500
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12500
$ws.Range("N43").Value = -12886
$ws.Range("M43").ClearContents()
# row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
# row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
# row 75
$ws.Range("H75").Value = 23779.75
$ws.Range("I75").Value = 22559.5
$ws.Range("J75").Value = 25000
$ws.Range("K75").Value = 22559.5
$ws.Range("L75").Value = 25000
$ws.Range("M75").Value = -21623.5
$ws.Range("N75").Value = -26872
# row 78
$ws.Range("H78").Value = 23779.75
$ws.Range("I78").Value = 22559.5
$ws.Range("J78").Value = 25000
$ws.Range("K78").Value = 67678.5
$ws.Range("L78").Value = 75000
$ws.Range("M78").Value = -62998.5
$ws.Range("N78").Value = -84360
# row 135
$ws.Range("H135").Value = 71207.25
$ws.Range("J135").Value = 71207.25
$ws.Range("L135").Value = 71207.25
$ws.Range("N135").Value = -81347.25
# row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# row 140
$ws.Range("H140").Value = 77214.5
$ws.Range("J140").Value = 77214.5
$ws.Range("L140").Value = 77214.5
$ws.Range("N140").Value = -87574.5
# row 141
$ws.Range("H141").Value = 87499.5
$ws.Range("J141").Value = 87499.5
$ws.Range("L141").Value = 87499.5
$ws.Range("N141").Value = -97859.5
